$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 22:20"

# --- Swap Costa Rica / Etiopia: Costa Rica now ranks above Etiopia ---
$ws.Range("A70").Value = "Costa Rica"
$ws.Range("A71").Value = "Etiopia"

# --- Refresh numeric stats (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4804123
$ws.Range("C4").Value = 39805
$ws.Range("D4").Value = 2375175
$ws.Range("E4").Value = 2270704
$ws.Range("G4").Value = 346
$ws.Range("H4").Value = 158244

# Row 8: Sudafrica
$ws.Range("B8").Value = 511485
$ws.Range("C8").Value = 8195
$ws.Range("D8").Value = 347227
$ws.Range("E8").Value = 155892
$ws.Range("G8").Value = 213
$ws.Range("H8").Value = 8366

# Row 70: now Costa Rica (updated figures)
$ws.Range("B70").Value = 18975
$ws.Range("C70").Value = 788
$ws.Range("D70").Value = 4585
$ws.Range("E70").Value = 14228
$ws.Range("G70").Value = 8
$ws.Range("H70").Value = 162

# Row 71: now Etiopia (figures carried over unchanged from the prior Costa Rica row slot)
$ws.Range("B71").Value = 18706
$ws.Range("C71").Value = 707
$ws.Range("D71").Value = 7601
$ws.Range("E71").Value = 10795
$ws.Range("G71").Value = 26
$ws.Range("H71").Value = 310

# Row 119: Suazilandia
$ws.Range("B119").Value = 2775
$ws.Range("C119").Value = 69
$ws.Range("E119").Value = 1518

# Row 131: Mozambique
$ws.Range("B131").Value = 1946
$ws.Range("C131").Value = 39
$ws.Range("D131").Value = 654
$ws.Range("E131").Value = 1279
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 13

# Row 160: Reunion
$ws.Range("B160").Value = 667
$ws.Range("C160").Value = 3
$ws.Range("E160").Value = 71

# Row 194: Belice
$ws.Range("B194").Value = 57
$ws.Range("C194").Value = 9
$ws.Range("E194").Value = 25
